$d = $word.ActiveDocument

# --- Fix 1: Title text on cover page ---
# Merge "CST 2540                Group Coursework " + "2" into one run,
# and "  G" + "9" into another run (this is what happens when the
# gramStart/gramEnd proofing marks collapse after retyping the title).
$d.Content.Find.Execute("CST 2540                Group Coursework 2  G9", $true, $false, $false, $false, $false, $true, 1, $false, "CST 2540                Group Coursework 2  G9", 2)

# --- Fix 2: Expand explanation in Q3 individual section ---
$d.Content.Find.Execute("returns to low.", $true, $false, $false, $false, $false, $true, 1, $false, "returns to low, and the machine remains in this state as long as X stays pressed after Z is done with its clock cycles.", 2)

$d.Content.Find.Execute("Upon X going back to low, the machine enters", $true, $false, $false, $false, $false, $true, 1, $false, "Upon X going back to low (which means the button is no longer being pressed), the machine enters", 2)
